$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped from
# 45186 (2023-09-17) to 45188 (2023-09-19) for every data row (rows 2-262).
$ws.Range("C2:C262").Value = 45188
